$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "41.867.92"
$ws.Cells.Item(2, 5).Value = "  -0.45%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.258.66"
$ws.Cells.Item(3, 5).Value = "  -0.60%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "304.95"
$ws.Cells.Item(5, 5).Value = "  -0.23%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "94.70"
$ws.Cells.Item(6, 5).Value = "  +1.58%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.525"
$ws.Cells.Item(7, 5).Value = "  -1.17%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.02%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.488"
$ws.Cells.Item(9, 5).Value = "  -0.22%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "34.89"
$ws.Cells.Item(10, 5).Value = "  +6.02%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0789"
$ws.Cells.Item(11, 5).Value = "  -1.87%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.44%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.63"
$ws.Cells.Item(13, 5).Value = "  -1.11%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.604.07"
$ws.Cells.Item(14, 5).Value = "  -0.72%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.37"
$ws.Cells.Item(15, 5).Value = "  -0.10%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.265.71"
$ws.Cells.Item(16, 5).Value = "  -0.33%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.789"
$ws.Cells.Item(17, 5).Value = "  +0.24%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "41.743.46"
$ws.Cells.Item(18, 5).Value = "  -0.36%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.32"
$ws.Cells.Item(19, 5).Value = "  -3.64%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0900"
$ws.Cells.Item(20, 5).Value = "  -2.19%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.96"
$ws.Cells.Item(21, 5).Value = "  -0.58%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "68.08"
$ws.Cells.Item(22, 5).Value = "  -0.29%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "236.82"
$ws.Cells.Item(23, 5).Value = "  -3.17%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.74%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.00"
$ws.Cells.Item(25, 5).Value = "  +0.12%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -1.50%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "23.61"
$ws.Cells.Item(27, 5).Value = "  -1.84%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "36.31"
$ws.Cells.Item(28, 5).Value = "  +3.22%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.25%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.43"
$ws.Cells.Item(30, 5).Value = "  -2.67%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "159.82"
$ws.Cells.Item(31, 5).Value = "  +0.06%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "5.21"
$ws.Cells.Item(32, 5).Value = "  -2.90%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.999"
$ws.Cells.Item(33, 5).Value = "  -0.03%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.15"
$ws.Cells.Item(34, 5).Value = "  +3.34%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0734"
$ws.Cells.Item(35, 5).Value = "  -1.69%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "16.96"
$ws.Cells.Item(36, 5).Value = "  -1.82%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.41%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.95%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.82"
$ws.Cells.Item(39, 5).Value = "  +0.87%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.114"
$ws.Cells.Item(40, 5).Value = "  -2.65%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.98"
$ws.Cells.Item(41, 5).Value = "  +0.30%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +3.21%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.962.91"
$ws.Cells.Item(43, 5).Value = "  -2.83%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0282"
$ws.Cells.Item(44, 5).Value = "  -0.52%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "18.63"
$ws.Cells.Item(45, 5).Value = "  -7.51%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.93"
$ws.Cells.Item(46, 5).Value = "  -0.18%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.84"
$ws.Cells.Item(47, 5).Value = "  -5.14%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "52.93"
$ws.Cells.Item(48, 5).Value = "  -0.88%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "72.73"
$ws.Cells.Item(49, 5).Value = "  +0.13%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.50"
$ws.Cells.Item(50, 5).Value = "  -1.82%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "90.87"
$ws.Cells.Item(51, 5).Value = "  -1.34%  "
